$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings used by A1 and D1) ---
$ws.Range("A1").Value = "Satisfaction from timely query resolution"
$ws.Range("D1").Value = "Overall Satisfaction"

# --- Data cell updates (train/test split & evaluation results) ---
$ws.Range("D4").Value = "Exceeded expectations"

$ws.Range("C6").Value = 4
$ws.Range("D6").Value = "Satisfied"

$ws.Range("B12").Value = 4
$ws.Range("D12").Value = "Satisfied"

$ws.Range("A15").Value = 3
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = "Satisfied"

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 4
$ws.Range("D17").Value = "Satisfied"

$ws.Range("D27").Value = "Exceeded expectations"

$ws.Range("A28").Value = 4
$ws.Range("D28").Value = "Satisfied"

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "Below expetations"

$ws.Range("A35").Value = 3

$ws.Range("B36").Value = 3
$ws.Range("C36").Value = 3

$ws.Range("A39").Value = 1
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = "Below expetations"

$ws.Range("A40").Value = 3
$ws.Range("C40").Value = 3

$ws.Range("A41").Value = 2
$ws.Range("B41").Value = 2
$ws.Range("C41").Value = 2

$ws.Range("A43").Value = 3
$ws.Range("B43").Value = 4
$ws.Range("C43").Value = 3
$ws.Range("D43").Value = "Satisfied"

$ws.Range("D44").Value = "Exceeded expectations"

$ws.Range("D45").Value = "Below expetations"

$ws.Range("D46").Value = "Exceeded expectations"

$ws.Range("C47").Value = 4

$ws.Range("D55").Value = "Below expetations"

$ws.Range("A61").Value = 3
$ws.Range("B61").Value = 3
$ws.Range("C61").Value = 3

$ws.Range("B62").Value = 4
$ws.Range("D62").Value = "Satisfied"

$ws.Range("D63").Value = "Below expetations"

$ws.Range("D71").Value = "Satisfied"

$ws.Range("D76").Value = "Exceeded expectations"

$ws.Range("D84").Value = "Exceeded expectations"

$ws.Range("D85").Value = "Exceeded expectations"

# --- View state: scroll so row 3 is the top-left visible row, select A27 ---
$ws.Range("A27").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
